# posted midterm 2 gradelines
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark additional homework / quiz items as graded (F column = 1.0)
$ws.Range("F9").Value = 1.0
$ws.Range("F10").Value = 1.0
$ws.Range("F11").Value = 1.0
$ws.Range("F21").Value = 1.0
$ws.Range("F22").Value = 1.0

# Midterm 2 total points now posted
$ws.Range("E32").Value = 85.0
$ws.Range("F32").Value = 1.0

# Midterm 2 gradeline scores (column L), rows 12-23 correspond to grades A .. F
$ws.Range("L12").Value = 72.0
$ws.Range("L13").Value = 67.0
$ws.Range("L14").Value = 62.0
$ws.Range("L15").Value = 57.0
$ws.Range("L16").Value = 52.0
$ws.Range("L17").Value = 46.0
$ws.Range("L18").Value = 40.0
$ws.Range("L19").Value = 34.0
$ws.Range("L20").Value = 29.0
$ws.Range("L21").Value = 25.0
$ws.Range("L22").Value = 20.0
$ws.Range("L23").Value = 0.0

# The O column ("Overall") formulas reference L12:L23, which just changed
# from blank to numeric. Re-assign each formula in place (no-op text-wise)
# so the dependency graph/recalc picks up the newly-populated precedents
# before the workbook is saved.
foreach ($r in 12..23) {
    $cell = $ws.Cells.Item($r, 15)
    $cell.Formula = $cell.Formula
}
